# "Generate Report for Handback"
#
# This localization-status workbook gets refreshed after a handback run:
#   - every "Ready for handoff" status cell becomes
#     "Handed back: in sync with en-US"
#   - the per-language "Latest Handback DateTime" timestamps advance to the
#     handback run's timestamps
#   - the stale "handback file is not the latest" Error Detail message is
#     cleared now that the handback is in sync
#   - the Status / Error Detail columns are resized to fit the new text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: per-language rollup status (E2 = zh-cn, F2 = de-de)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# zh-cn detail sheet
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-09-05 18:55:53"
$wsZhCn.Range("P2").Value = ""

# de-de detail sheet
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-09-05 18:56:01"
$wsDeDe.Range("P2").Value = ""

# Resize the Status columns (now holding the longer status text) and the
# Error Detail columns (now empty) to match the refreshed report layout.
$wsOverview.Columns("E").ColumnWidth = 29.144371396019398
$wsOverview.Columns("F").ColumnWidth = 29.144371396019398

$wsZhCn.Columns("C").ColumnWidth = 29.144371396019398
$wsZhCn.Columns("P").ColumnWidth = 12.913719540550602

$wsDeDe.Columns("C").ColumnWidth = 29.144371396019398
$wsDeDe.Columns("P").ColumnWidth = 12.913719540550602
